$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 13.843044
$ws.Range("H2").Value = 41.529132
$ws.Range("I2").Value = 0.1139869403149299
$ws.Range("J2").Value = 0.1139869403149299
$ws.Range("M2").Value = 133.3951123333333
$ws.Range("N2").Value = 400.185337
$ws.Range("O2").Value = 0.8984588679103155
$ws.Range("P2").Value = 0.8984588679103156
$ws.Range("Q2").Value = 1846.594409415276
$ws.Range("R2").Value = 16619.34968473748
$ws.Range("S2").Value = 0.1024125773519126
$ws.Range("T2").Value = 0.1024125773519126
$ws.Range("G3").Value = 13.843044
$ws.Range("H3").Value = 41.529132
$ws.Range("I3").Value = 0.1139869403149299
$ws.Range("J3").Value = 0.1139869403149299
$ws.Range("M3").Value = 2.340788333333334
$ws.Range("N3").Value = 7.022365000000001
$ws.Range("O3").Value = 0.01576596023045448
$ws.Range("P3").Value = 0.01576596023045448
$ws.Range("Q3").Value = 32.40363589302
$ws.Range("R3").Value = 291.63272303718
$ws.Range("S3").Value = 0.001797113567796372
$ws.Range("T3").Value = 0.001797113567796373
$ws.Range("G4").Value = 13.843044
$ws.Range("H4").Value = 41.529132
$ws.Range("I4").Value = 0.1139869403149299
$ws.Range("J4").Value = 0.1139869403149299
$ws.Range("M4").Value = 12.735128
$ws.Range("N4").Value = 38.205384
$ws.Range("O4").Value = 0.08577517185923002
$ws.Range("P4").Value = 0.08577517185923003
$ws.Range("Q4").Value = 176.292937249632
$ws.Range("R4").Value = 1586.636435246688
$ws.Range("S4").Value = 0.009777249395220903
$ws.Range("T4").Value = 0.009777249395220905
$ws.Range("I5").Value = 0.7803892412315415
$ws.Range("J5").Value = 0.7803892412315415
$ws.Range("M5").Value = 133.3951123333333
$ws.Range("N5").Value = 400.185337
$ws.Range("O5").Value = 0.8984588679103155
$ws.Range("P5").Value = 0.8984588679103156
$ws.Range("Q5").Value = 12642.34662360917
$ws.Range("R5").Value = 113781.1196124826
$ws.Range("S5").Value = 0.7011476342062809
$ws.Range("T5").Value = 0.701147634206281
$ws.Range("I6").Value = 0.7803892412315415
$ws.Range("J6").Value = 0.7803892412315415
$ws.Range("M6").Value = 2.340788333333334
$ws.Range("N6").Value = 7.022365000000001
$ws.Range("O6").Value = 0.01576596023045448
$ws.Range("P6").Value = 0.01576596023045448
$ws.Range("S6").Value = 0.01230358574153103
$ws.Range("T6").Value = 0.01230358574153103
$ws.Range("I7").Value = 0.7803892412315415
$ws.Range("J7").Value = 0.7803892412315415
$ws.Range("M7").Value = 12.735128
$ws.Range("N7").Value = 38.205384
$ws.Range("O7").Value = 0.08577517185923002
$ws.Range("P7").Value = 0.08577517185923003
$ws.Range("Q7").Value = 1206.955034976936
$ws.Range("R7").Value = 10862.59531479242
$ws.Range("S7").Value = 0.06693802128372958
$ws.Range("T7").Value = 0.0669380212837296
$ws.Range("G8").Value = 12.827392
$ws.Range("H8").Value = 38.482176
$ws.Range("I8").Value = 0.1056238184535286
$ws.Range("J8").Value = 0.1056238184535286
$ws.Range("M8").Value = 133.3951123333333
$ws.Range("N8").Value = 400.185337
$ws.Range("O8").Value = 0.8984588679103155
$ws.Range("P8").Value = 0.8984588679103156
$ws.Range("Q8").Value = 1711.111396783702
$ws.Range("R8").Value = 15400.00257105331
$ws.Range("S8").Value = 0.09489865635212204
$ws.Range("T8").Value = 0.09489865635212204
$ws.Range("G9").Value = 12.827392
$ws.Range("H9").Value = 38.482176
$ws.Range("I9").Value = 0.1056238184535286
$ws.Range("J9").Value = 0.1056238184535286
$ws.Range("M9").Value = 2.340788333333334
$ws.Range("N9").Value = 7.022365000000001
$ws.Range("O9").Value = 0.01576596023045448
$ws.Range("P9").Value = 0.01576596023045448
$ws.Range("Q9").Value = 30.02620954069334
$ws.Range("R9").Value = 270.23588586624
$ws.Range("S9").Value = 0.001665260921127076
$ws.Range("T9").Value = 0.001665260921127077
$ws.Range("G10").Value = 12.827392
$ws.Range("H10").Value = 38.482176
$ws.Range("I10").Value = 0.1056238184535286
$ws.Range("J10").Value = 0.1056238184535286
$ws.Range("M10").Value = 12.735128
$ws.Range("N10").Value = 38.205384
$ws.Range("O10").Value = 0.08577517185923002
$ws.Range("P10").Value = 0.08577517185923003
$ws.Range("Q10").Value = 163.358479026176
$ws.Range("R10").Value = 1470.226311235584
$ws.Range("S10").Value = 0.00905990118027953
$ws.Range("T10").Value = 0.009059901180279531
